$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 332.5
$ws.Range("H32").Value = 1099
$ws.Range("I32").Value = 995
$ws.Range("J32").Value = 1116.3334
$ws.Range("K32").Value = 995
$ws.Range("L32").Value = 1116.3334
$ws.Range("M32").Value = -669
$ws.Range("N32").Value = -1768.3334
$ws.Range("H76").Value = 3387.5
$ws.Range("I76").Value = 3414.2856
$ws.Range("K76").Value = 3414.2856
$ws.Range("M76").Value = -3099.2856
$ws.Range("H79").Value = 3387.5
$ws.Range("I79").Value = 3414.2856
$ws.Range("K79").Value = 3414.2856
$ws.Range("M79").Value = -2322.2856
$ws.Range("H132").Value = 15433752
$ws.Range("I132").Value = 1884042.4
$ws.Range("J132").Value = 76928590
$ws.Range("K132").Value = 5652127.199999999
$ws.Range("L132").Value = 230785770
$ws.Range("M132").Value = -5649597.199999999
$ws.Range("N132").Value = -230790830
$ws.Range("H133").Value = 34950
$ws.Range("J133").Value = 34950
$ws.Range("L133").Value = 34950
$ws.Range("N133").Value = -45070
$ws.Range("H138").Value = 3924530
$ws.Range("I138").Value = 1374.1
$ws.Range("J138").Value = 6455598.5
$ws.Range("K138").Value = 4122.299999999999
$ws.Range("L138").Value = 19366795.5
$ws.Range("M138").Value = 1017.700000000001
$ws.Range("N138").Value = -19377075.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 400
$ws.Range("I16").Value = 400
$ws.Range("K16").Value = 400
$ws.Range("M16").Value = -113
$ws.Range("H32").Value = 4793.7534
$ws.Range("I32").Value = 3327.75
$ws.Range("K32").Value = 3327.75
$ws.Range("M32").Value = -3040.75
$ws.Range("H50").Value = 1104
$ws.Range("I50").Value = 684.5
$ws.Range("J50").Value = 1663.3334
$ws.Range("K50").Value = 684.5
$ws.Range("L50").Value = 1663.3334
$ws.Range("M50").Value = 29.5
$ws.Range("N50").Value = -3091.3334
$ws.Range("H53").Value = 20492.375
$ws.Range("I53").Value = 7969.5
$ws.Range("J53").Value = 24666.666
$ws.Range("K53").Value = 7969.5
$ws.Range("L53").Value = 24666.666
$ws.Range("M53").Value = -7287.5
$ws.Range("N53").Value = -26030.666
$ws.Range("H61").Value = 1384.3889
$ws.Range("I61").Value = 1419.9412
$ws.Range("J61").Value = 780
$ws.Range("K61").Value = 1419.9412
$ws.Range("L61").Value = 780
$ws.Range("M61").Value = -1207.9412
$ws.Range("N61").Value = -1204
$ws.Range("H63").Value = 3027.818
$ws.Range("I63").Value = 2116.6667
$ws.Range("J63").Value = 4121.2
$ws.Range("K63").Value = 2116.6667
$ws.Range("L63").Value = 4121.2
$ws.Range("M63").Value = -1430.6667
$ws.Range("N63").Value = -5493.2
$ws.Range("H66").Value = 3027.818
$ws.Range("I66").Value = 2116.6667
$ws.Range("J66").Value = 4121.2
$ws.Range("K66").Value = 10583.3335
$ws.Range("L66").Value = 20606
$ws.Range("M66").Value = -7151.333500000001
$ws.Range("N66").Value = -27470
$ws.Range("H88").Value = 2333.3333
$ws.Range("I88").Value = 2333.3333
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 2333.3333
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -1927.3333
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 2333.3333
$ws.Range("I91").Value = 2333.3333
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 2333.3333
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -929.3332999999998
$ws.Range("N91").ClearContents()
$ws.Range("H136").Value = 1384.3889
$ws.Range("I136").Value = 1419.9412
$ws.Range("J136").Value = 780
$ws.Range("K136").Value = 4259.8236
$ws.Range("L136").Value = 2340
$ws.Range("M136").Value = -1709.8236
$ws.Range("N136").Value = -7440

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 29169.5
$ws.Range("J55").Value = 29169.5
$ws.Range("L55").Value = 29169.5
$ws.Range("N55").Value = -29715.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 6994
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 6994
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 6994
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -7334
$ws.Range("H31").Value = 19232356
$ws.Range("I31").Value = 34483652
$ws.Range("J31").Value = 2462.6086
$ws.Range("K31").Value = 34483652
$ws.Range("L31").Value = 2462.6086
$ws.Range("M31").Value = -34483357
$ws.Range("N31").Value = -3052.6086
$ws.Range("H34").Value = 19232356
$ws.Range("I34").Value = 34483652
$ws.Range("J34").Value = 2462.6086
$ws.Range("K34").Value = 34483652
$ws.Range("L34").Value = 2462.6086
$ws.Range("M34").Value = -34483450
$ws.Range("N34").Value = -2866.6086
$ws.Range("H99").Value = 3476.641
$ws.Range("I99").Value = 3378.5186
$ws.Range("J99").Value = 3697.4167
$ws.Range("K99").Value = 3378.5186
$ws.Range("L99").Value = 3697.4167
$ws.Range("M99").Value = -1880.5186
$ws.Range("N99").Value = -6693.4167
$ws.Range("H126").Value = 3476.641
$ws.Range("I126").Value = 3378.5186
$ws.Range("J126").Value = 3697.4167
$ws.Range("K126").Value = 10135.5558
$ws.Range("L126").Value = 11092.2501
$ws.Range("M126").Value = -7665.5558
$ws.Range("N126").Value = -16032.2501
$ws.Range("H132").Value = 3639
$ws.Range("I132").Value = 3460.5386
$ws.Range("J132").Value = 3970.4285
$ws.Range("K132").Value = 10381.6158
$ws.Range("L132").Value = 11911.2855
$ws.Range("M132").Value = -7851.6158
$ws.Range("N132").Value = -16971.2855
$ws.Range("H134").Value = 47828724
$ws.Range("I134").Value = 5558515
$ws.Range("J134").Value = 200001470
$ws.Range("K134").Value = 16675545
$ws.Range("L134").Value = 600004410
$ws.Range("M134").Value = -16673010
$ws.Range("N134").Value = -600009480

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 189.48276
$ws.Range("I12").Value = 261.22223
$ws.Range("J12").Value = 157.2
$ws.Range("K12").Value = 783.66669
$ws.Range("L12").Value = 471.6
$ws.Range("M12").Value = -610.66669
$ws.Range("N12").Value = -817.5999999999999
$ws.Range("H80").Value = 2430
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2430
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 7290
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -9162
$ws.Range("H83").Value = 2430
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2430
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 21870
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -31230
$ws.Range("H122").Value = 957.12
$ws.Range("I122").Value = 583.6429000000001
$ws.Range("J122").Value = 1432.4546
$ws.Range("K122").Value = 5252.7861
$ws.Range("L122").Value = 12892.0914
$ws.Range("M122").Value = -2802.7861
$ws.Range("N122").Value = -17792.0914

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5133.5
$ws.Range("I70").Value = 4969.8696
$ws.Range("J70").Value = 5354.8823
$ws.Range("K70").Value = 4969.8696
$ws.Range("L70").Value = 5354.8823
$ws.Range("M70").Value = -4699.8696
$ws.Range("N70").Value = -5894.8823
$ws.Range("H73").Value = 5133.5
$ws.Range("I73").Value = 4969.8696
$ws.Range("J73").Value = 5354.8823
$ws.Range("K73").Value = 4969.8696
$ws.Range("L73").Value = 5354.8823
$ws.Range("M73").Value = -4033.8696
$ws.Range("N73").Value = -7226.8823
$ws.Range("H80").Value = 2467.9285
$ws.Range("I80").Value = 2420.9167
$ws.Range("J80").Value = 2750
$ws.Range("K80").Value = 2420.9167
$ws.Range("L80").Value = 2750
$ws.Range("M80").Value = -1422.9167
$ws.Range("N80").Value = -4746
$ws.Range("H83").Value = 2467.9285
$ws.Range("I83").Value = 2420.9167
$ws.Range("J83").Value = 2750
$ws.Range("K83").Value = 12104.5835
$ws.Range("L83").Value = 13750
$ws.Range("M83").Value = -7112.583500000001
$ws.Range("N83").Value = -23734

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 29500
$ws.Range("J64").Value = 29500
$ws.Range("L64").Value = 29500
$ws.Range("N64").Value = -29996
$ws.Range("H67").Value = 29500
$ws.Range("J67").Value = 29500
$ws.Range("L67").Value = 29500
$ws.Range("N67").Value = -31216
$ws.Range("H93").Value = 30194.5
$ws.Range("J93").Value = 30194.5
$ws.Range("L93").Value = 30194.5
$ws.Range("N93").Value = -35186.5

